$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wells table")

# --- Insert three new columns (mirrors Excel's "Insert Sheet Columns") ---
# 1) New column before old "E" (S velocity) -> becomes "Sonic"
$dWidth = $ws.Range("D1").ColumnWidth
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").ColumnWidth = $dWidth
# 2) New column before old "F" (Gamma ray, now shifted to G) -> becomes "Shear Sonic"
$fWidth = $ws.Range("F1").ColumnWidth
$ws.Range("G1").EntireColumn.Insert()
$ws.Range("G1").ColumnWidth = $fWidth
# 3) New column before old "M" (Depth, now shifted to O) -> becomes "Caliper"
$ws.Range("O1").EntireColumn.Insert()

# --- New header cells for the inserted columns ---
$ws.Range("E2").Value = "Sonic"

# --- Update "Use" column values (A) for several rows ---
$ws.Range("A4").Value = "No"
$ws.Range("A5").Value = "No"
$ws.Range("A7").Value = "No"
$ws.Range("A9").Value = "No"

# --- Fill in new alias info for well "Well_L" (row 9), and the remaining
# header cells, written in the order that matches the author's original
# shared-string insertion order.
$ws.Range("E9").Value = "DT"
$ws.Range("G2").Value = "Shear Sonic"
$ws.Range("O9").Value = "CALD, CALS"
$ws.Range("H9").Value = "GRD, GRS"
$ws.Range("N9").Value = "ILD, ILM, LL8"
$ws.Range("I9").Value = "RHOB"
$ws.Range("J9").Value = "NPHILS, NPHISS"
$ws.Range("A8").Value = "NO"

$ws.Range("O2").Value = "Caliper"

# --- Move the active-cell selection (cosmetic, matches author's saved view) ---
$ws.Range("A10").Select()
